$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (paragraph 2, right after the
#    Heading1 title). Deleting the paragraph's Range (including its mark)
#    merges paragraph 1 directly onto what was paragraph 3.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Turn the final paragraph (the italic AI image-prompt paragraph) into
#    the new bold "Play Derby Dash..." paragraph, and put the old
#    "Meta description" wording (minus the "Meta description" label) into
#    a freshly appended final paragraph, keeping the italic look.
$d.Paragraphs.Add() | Out-Null

$count = $d.Paragraphs.Count
$headingPara = $d.Paragraphs($count - 1)
$descPara = $d.Paragraphs($count)

$headingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingRange.Text = "Play Derby Dash for Free - Review and Ratings"
$headingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingRange.Font.Italic = $false
$headingRange.Font.Bold = $true

$descRange = $d.Range($descPara.Range.Start, $descPara.Range.Start)
$descRange.InsertAfter("Find out all about Derby Dash, its features, winning chances and mobile compatibility, and play it for free.")
